$wb = $excel.ActiveWorkbook

# 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1644
$ws1.Range("F8").Value = 2144
$ws1.Range("F15").Value = 2079
$ws1.Range("F18").Value = 2572
$ws1.Range("F19").Value = 33
$ws1.Range("F21").Value = 310
$ws1.Range("F27").Value = 4444

# 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 55
$ws2.Range("F14").Value = 305

# 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1765
$ws3.Range("F7").Value = 460
$ws3.Range("F8").Value = 77

# 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1765
$ws4.Range("F5").Value = 460
$ws4.Range("F11").Value = 1644
$ws4.Range("F17").Value = 2144
$ws4.Range("F23").Value = 55
$ws4.Range("F27").Value = 305
$ws4.Range("F29").Value = 2079
$ws4.Range("F34").Value = 2572
$ws4.Range("F36").Value = 33
$ws4.Range("F46").Value = 4444
